{"js": "const replacements = [\n  [\"19\u00d747=893\", \"37\u00d747=1739\"],\n  [\"57\u00d752=2964\", \"69\u00d737=2553\"],\n  [\"25\u00d751=1275\", \"46\u00d729=1334\"],\n  [\"55\u00d795=5225\", \"76\u00d730=2280\"],\n  [\"24\u00d766=1584\", \"71\u00d781=5751\"],\n  [\"28\u00d787=2436\", \"46\u00d765=2990\"],\n  [\"96\u00d798=9408\", \"58\u00d777=4466\"],\n  [\"59\u00d776=4484\", \"95\u00d722=2090\"],\n  [\"92\u00d757=5244\", \"35\u00d741=1435\"],\n  [\"69\u00d758=4002\", \"90\u00d778=7020\"],\n  [\"95\u00d771=6745\", \"12\u00d794=1128\"],\n  [\"24\u00d761=1464\", \"99\u00d776=7524\"],\n  [\"35\u00d779=2765\", \"26\u00d734=884\"],\n  [\"69\u00d718=1242\", \"53\u00d756=2968\"],\n  [\"18\u00d767=1206\", \"26\u00d773=1898\"],\n  [\"23\u00d716=368\", \"39\u00d785=3315\"],\n  [\"93\u00d742=3906\", \"38\u00d726=988\"],\n  [\"38\u00d730=1140\", \"15\u00d747=705\"],\n  [\"35\u00d785=2975\", \"72\u00d793=6696\"],\n  [\"36\u00d791=3276\", \"14\u00d735=490\"],\n  [\"63\u00d717=1071\", \"75\u00d741=3075\"],\n  [\"34\u00d764=2176\", \"76\u00d753=4028\"],\n  [\"83\u00d761=5063\", \"45\u00d750=2250\"],\n  [\"92\u00d772=6624\", \"40\u00d755=2200\"],\n  [\"81\u00d726=2106\", \"13\u00d798=1274\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"19\u00d747=893\", \"37\u00d747=1739\"),\n    @(\"57\u00d752=2964\", \"69\u00d737=2553\"),\n    @(\"25\u00d751=1275\", \"46\u00d729=1334\"),\n    @(\"55\u00d795=5225\", \"76\u00d730=2280\"),\n    @(\"24\u00d766=1584\", \"71\u00d781=5751\"),\n    @(\"28\u00d787=2436\", \"46\u00d765=2990\"),\n    @(\"96\u00d798=9408\", \"58\u00d777=4466\"),\n    @(\"59\u00d776=4484\", \"95\u00d722=2090\"),\n    @(\"92\u00d757=5244\", \"35\u00d741=1435\"),\n    @(\"69\u00d758=4002\", \"90\u00d778=7020\"),\n    @(\"95\u00d771=6745\", \"12\u00d794=1128\"),\n    @(\"24\u00d761=1464\", \"99\u00d776=7524\"),\n    @(\"35\u00d779=2765\", \"26\u00d734=884\"),\n    @(\"69\u00d718=1242\", \"53\u00d756=2968\"),\n    @(\"18\u00d767=1206\", \"26\u00d773=1898\"),\n    @(\"23\u00d716=368\", \"39\u00d785=3315\"),\n    @(\"93\u00d742=3906\", \"38\u00d726=988\"),\n    @(\"38\u00d730=1140\", \"15\u00d747=705\"),\n    @(\"35\u00d785=2975\", \"72\u00d793=6696\"),\n    @(\"36\u00d791=3276\", \"14\u00d735=490\"),\n    @(\"63\u00d717=1071\", \"75\u00d741=3075\"),\n    @(\"34\u00d764=2176\", \"76\u00d753=4028\"),\n    @(\"83\u00d761=5063\", \"45\u00d750=2250\"),\n    @(\"92\u00d772=6624\", \"40\u00d755=2200\"),\n    @(\"81\u00d726=2106\", \"13\u00d798=1274\"),\n)\n\n$wdReplaceAll = 2\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, $wdReplaceAll)\n}\n"}
